# "add some new commands to excel" - insert a new Git command row
# (git commit -am ...) right after the existing "git commit -m" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 6 down by one to make room for the new command.
$ws.Rows.Item(6).Insert()

# Fill B6 before A6 so the shared-string table gets the new strings in the
# same order as the reference workbook (B's text becomes index 39, A's
# text becomes index 40).
$ws.Range("B6").Value = "am means git commit all of the files that have been changed--"
$ws.Range("A6").Value = 'git commit -am "text abou commit" '

# Inserting the row duplicated the tall (25.5pt) row height from the rows
# that used to need two lines of text; those rows are single-line again
# now that the new row is its own thing, so auto-fit them back down.
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).AutoFit()

# Leave the selection on the newly added cell.
[void]$ws.Range("A6").Select()

Write-Output "done"
